$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = -0.4458934228143729
    "C2"  = -0.1134040234363881
    "D2"  = 0.6219785646361177
    "B3"  = -0.04237791528613299
    "C3"  = 0.3807011644083477
    "B4"  = 0.7181631419195285
    "B5"  = 1.049089967442578
    "C5"  = 0.8524940945699551
    "D5"  = -0.1034953623938322
    "E5"  = 0.0926409719578045
    "B6"  = 1.499698792741737
    "C6"  = -0.0441541421789331
    "D6"  = -0.3287619841940873
    "B7"  = 0.8297131584907731
    "C7"  = -0.1093333813988623
    "B8"  = 0.1620838256790951
    "B9"  = 0.3211787666311243
    "C9"  = 0.112909106075748
    "D9"  = 0.03134352168291979
    "E9"  = -0.06584948318657535
    "B10" = 0.2286550815363005
    "C10" = 0.125655419861192
    "D10" = -0.1481903203942864
    "B11" = 0.4490319412376911
    "C11" = -0.2127353874755186
    "B12" = 0.1181141524322671
    "B13" = -0.0001442643046098147
    "C13" = -0.008000482844254697
    "D13" = 0.2025442913845202
    "B14" = 0.1064660537556523
    "C14" = 0.1601369838950918
    "B15" = 0.04651757380517277
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
